# Generate Report for Handoff
# This script updates the localization-status workbook from a "handed back" snapshot
# to a "ready for handoff" snapshot: new GUID-named files, new status text, new
# timestamps, and removal of the "Latest Target File" / "Latest Handback File"
# columns' data (the handback has not happened yet).

$wb = $excel.ActiveWorkbook

$newMdA      = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
$newMdB      = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md"
$newStatus   = "Ready for handoff"
$newDateD    = "2016-03-24 03:14:40"
$newHandoffE = "2016-03-24 03:14:36"
$newHandback = "0001-01-01 00:00:00"
$newZhXlf    = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$newDeXlf    = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Cells.Item(2,1).Value = $newMdA
$wsOverview.Cells.Item(2,2).Value = $newStatus
$wsOverview.Cells.Item(2,3).Value = $newStatus
$wsOverview.Cells.Item(2,4).Value = $newDateD

$wsOverview.Cells.Item(3,1).Value = $newMdB
$wsOverview.Cells.Item(3,2).Value = $newStatus
$wsOverview.Cells.Item(3,3).Value = $newStatus
$wsOverview.Cells.Item(3,4).Value = $newDateD

$ovA2Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/8d98bf3f-aa28-4fc1-808e-724969253c35.md"
$ovA3Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/f8955185-2377-4935-980b-9748bd6ee4d5.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $ovA2Address, "", "", $newMdA)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $ovA3Address, "", "", $newMdB)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$wsZh.Cells.Item(2,1).Value = $newMdA
$wsZh.Cells.Item(2,3).Value = $newStatus
$wsZh.Cells.Item(2,4).Value = $newZhXlf
$wsZh.Cells.Item(2,5).Value = $newHandoffE
$wsZh.Cells.Item(2,6).Clear()
$wsZh.Cells.Item(2,7).Clear()
$wsZh.Cells.Item(2,8).Value = $newHandback

$wsZh.Cells.Item(3,1).Value = $newMdB
$wsZh.Cells.Item(3,3).Value = $newStatus
$wsZh.Cells.Item(3,4).Value = $newZhXlf
$wsZh.Cells.Item(3,5).Value = $newHandoffE
$wsZh.Cells.Item(3,6).Clear()
$wsZh.Cells.Item(3,7).Clear()
$wsZh.Cells.Item(3,8).Value = $newHandback

$zhA2Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/8d98bf3f-aa28-4fc1-808e-724969253c35.md"
$zhD2Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/008040ab84ec282210d48190becc24751080f2cd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8d98bf3f-aa28-4fc1-808e-724969253c35.c3231ca7600b20953891f06c61acff0fb05dea4c.zh-cn.xlf"
$zhA3Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/f8955185-2377-4935-980b-9748bd6ee4d5.md"
$zhD3Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/008040ab84ec282210d48190becc24751080f2cd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f8955185-2377-4935-980b-9748bd6ee4d5.ea6c6c93cdd8497d730f89bf41785b4448341f59.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Address, "", "", $newMdA)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhD2Address, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Address, "", "", $newMdB)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhD3Address, "", "", $newZhXlf)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$wsDe.Cells.Item(2,1).Value = $newMdA
$wsDe.Cells.Item(2,3).Value = $newStatus
$wsDe.Cells.Item(2,4).Value = $newDeXlf
$wsDe.Cells.Item(2,5).Value = $newDateD
$wsDe.Cells.Item(2,6).Clear()
$wsDe.Cells.Item(2,7).Clear()
$wsDe.Cells.Item(2,8).Value = $newHandback

$wsDe.Cells.Item(3,1).Value = $newMdB
$wsDe.Cells.Item(3,3).Value = $newStatus
$wsDe.Cells.Item(3,4).Value = $newDeXlf
$wsDe.Cells.Item(3,5).Value = $newDateD
$wsDe.Cells.Item(3,6).Clear()
$wsDe.Cells.Item(3,7).Clear()
$wsDe.Cells.Item(3,8).Value = $newHandback

$deA2Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/8d98bf3f-aa28-4fc1-808e-724969253c35.md"
$deD2Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89a58e244a18752b927a1a8a7810e1acdeef3d6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8d98bf3f-aa28-4fc1-808e-724969253c35.c3231ca7600b20953891f06c61acff0fb05dea4c.de-de.xlf"
$deA3Address = "https://github.com/OpenLocalizationTest/oltest/blob/508564eb20dde7797b340609612ea1311afc52c6/e2e/f8955185-2377-4935-980b-9748bd6ee4d5.md"
$deD3Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89a58e244a18752b927a1a8a7810e1acdeef3d6f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f8955185-2377-4935-980b-9748bd6ee4d5.ea6c6c93cdd8497d730f89bf41785b4448341f59.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Address, "", "", $newMdA)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deD2Address, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Address, "", "", $newMdB)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deD3Address, "", "", $newDeXlf)

Write-Output "Report regenerated for handoff."
